$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -6
$ws.Range("F5").Value = 1
$ws.Range("F16").Value = 2
$ws.Range("F17").Value = -3
$ws.Range("F18").Value = -1
